$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last four retired test-scenario rows (102-105)
$ws.Range("A102:C105").EntireRow.Delete()

# Rewrite rows 80-101 with the updated set of test scenarios
$ws.Cells.Item(80, 1).Value = "Schools  Position Catagories Delete Functionality"
$ws.Cells.Item(80, 2).Value = "FAILED"
$ws.Cells.Item(80, 3).Value = "chrome"
$ws.Cells.Item(81, 1).Value = "Schools  Position Catagories Delete Functionality"
$ws.Cells.Item(81, 2).Value = "FAILED"
$ws.Cells.Item(81, 3).Value = "chrome"
$ws.Cells.Item(82, 1).Value = "Human Resources  Position Catagories Add Functionality"
$ws.Cells.Item(82, 2).Value = "FAILED"
$ws.Cells.Item(82, 3).Value = "edge"
$ws.Cells.Item(83, 1).Value = "Scholls  Position Catagories Add Functionality"
$ws.Cells.Item(83, 2).Value = "FAILED"
$ws.Cells.Item(83, 3).Value = "edge"
$ws.Cells.Item(84, 1).Value = "Human Resources  Position Catagories Add Functionality"
$ws.Cells.Item(84, 2).Value = "PASSED"
$ws.Cells.Item(84, 3).Value = "chrome"
$ws.Cells.Item(85, 1).Value = "Scholls  Position Catagories Add Functionality"
$ws.Cells.Item(85, 2).Value = "FAILED"
$ws.Cells.Item(85, 3).Value = "chrome"
$ws.Cells.Item(86, 1).Value = "Human Resources  Position Catagories Add Functionality"
$ws.Cells.Item(86, 2).Value = "FAILED"
$ws.Cells.Item(86, 3).Value = "chrome"
$ws.Cells.Item(87, 1).Value = "Scholls  Position Catagories Add Functionality"
$ws.Cells.Item(87, 2).Value = "FAILED"
$ws.Cells.Item(87, 3).Value = "chrome"
$ws.Cells.Item(88, 1).Value = "Human Resources  Attestations Add Functionality"
$ws.Cells.Item(88, 2).Value = "FAILED"
$ws.Cells.Item(88, 3).Value = "chrome"
$ws.Cells.Item(89, 1).Value = "Human Resources  Attestations  Negative Add Functionality"
$ws.Cells.Item(89, 2).Value = "PASSED"
$ws.Cells.Item(89, 3).Value = "chrome"
$ws.Cells.Item(90, 1).Value = "Human Resources  Attestations  Edit Functionality"
$ws.Cells.Item(90, 2).Value = "PASSED"
$ws.Cells.Item(90, 3).Value = "chrome"
$ws.Cells.Item(91, 1).Value = "Human Resources  Attestations Delete Functionality"
$ws.Cells.Item(91, 2).Value = "FAILED"
$ws.Cells.Item(91, 3).Value = "chrome"
$ws.Cells.Item(92, 1).Value = "Human Resources  Attestations Negative Delete Functionality"
$ws.Cells.Item(92, 2).Value = "FAILED"
$ws.Cells.Item(92, 3).Value = "chrome"
$ws.Cells.Item(93, 1).Value = "Scholls  Position Catagories Add Functionality"
$ws.Cells.Item(93, 2).Value = "FAILED"
$ws.Cells.Item(93, 3).Value = "chrome"
$ws.Cells.Item(94, 1).Value = "Scholls  Position Catagories Negative Add Functionality"
$ws.Cells.Item(94, 2).Value = "FAILED"
$ws.Cells.Item(94, 3).Value = "chrome"
$ws.Cells.Item(95, 1).Value = "Schools  Position Catagories Edit Functionality"
$ws.Cells.Item(95, 2).Value = "FAILED"
$ws.Cells.Item(95, 3).Value = "chrome"
$ws.Cells.Item(96, 1).Value = "Scholls  Position Catagories Add Functionality"
$ws.Cells.Item(96, 2).Value = "FAILED"
$ws.Cells.Item(96, 3).Value = "chrome"
$ws.Cells.Item(97, 1).Value = "Scholls  Position Catagories Add Functionality"
$ws.Cells.Item(97, 2).Value = "FAILED"
$ws.Cells.Item(97, 3).Value = "chrome"
$ws.Cells.Item(98, 1).Value = "Scholls  Position Catagories Add Functionality"
$ws.Cells.Item(98, 2).Value = "FAILED"
$ws.Cells.Item(98, 3).Value = "chrome"
$ws.Cells.Item(99, 1).Value = "Scholls  Position Catagories Add Functionality"
$ws.Cells.Item(99, 2).Value = "PASSED"
$ws.Cells.Item(99, 3).Value = "chrome"
$ws.Cells.Item(100, 1).Value = "Scholls  Position Catagories Negative Add Functionality"
$ws.Cells.Item(100, 2).Value = "FAILED"
$ws.Cells.Item(100, 3).Value = "chrome"
$ws.Cells.Item(101, 1).Value = "Scholls  Position Catagories Negative Add Functionality"
$ws.Cells.Item(101, 2).Value = "PASSED"
$ws.Cells.Item(101, 3).Value = "chrome"
